$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individual Costs")

$ws.Cells.Item(2, 3).Value = 242.9301789682836
$ws.Cells.Item(3, 3).Value = 241.7533522765146
$ws.Cells.Item(4, 3).Value = 244.4072642249967
$ws.Cells.Item(5, 3).Value = 230.6829339348994
$ws.Cells.Item(6, 3).Value = 219.1516994869718
$ws.Cells.Item(7, 3).Value = 211.218754640267
$ws.Cells.Item(8, 3).Value = 203.3390079515765
$ws.Cells.Item(9, 3).Value = 201.3797721673737
$ws.Cells.Item(10, 3).Value = 245.4722048135623
$ws.Cells.Item(11, 3).Value = 254.8313851478964
$ws.Cells.Item(12, 3).Value = 260.2003103654174
$ws.Cells.Item(13, 3).Value = 259.9537246242008
$ws.Cells.Item(14, 3).Value = 252.9227432672594
$ws.Cells.Item(15, 3).Value = 267.3152168258129
$ws.Cells.Item(16, 3).Value = 300.3356876444313
$ws.Cells.Item(17, 3).Value = 305.7807857621891
$ws.Cells.Item(18, 3).Value = 258.7384440859091
$ws.Cells.Item(19, 3).Value = 237.7027714931847
$ws.Cells.Item(20, 3).Value = 281.3168672212262
$ws.Cells.Item(21, 3).Value = 256.801896043078
$ws.Cells.Item(22, 3).Value = 238.1904045372331
$ws.Cells.Item(23, 3).Value = 229.3050071278335
$ws.Cells.Item(24, 3).Value = 203.1075069825606
$ws.Cells.Item(25, 3).Value = 217.5450755808948
